$wb = $excel.ActiveWorkbook

$tor = $wb.Worksheets.Item("TOR")
$nontor = $wb.Worksheets.Item("Non-TOR")

# --- 1. Snapshot the two sheets' current (pre-edit) data into two new
#        "_Initial" sheets appended at the end of the workbook, as plain
#        values (no column-width formatting carried over). ---

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$torInitial = $wb.Worksheets.Add($null, $lastSheet)
$torInitial.Name = "TOR_Initial"
$tor.Range("A1:AB6").Copy() | Out-Null
$torInitial.Range("A1").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false
$torInitial.Range("A1:AB6").Select() | Out-Null

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$nontorInitial = $wb.Worksheets.Add($null, $lastSheet)
$nontorInitial.Name = "Non-TOR_Initial"
$nontor.Range("A1:AB6").Copy() | Out-Null
$nontorInitial.Range("A1").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false
$nontorInitial.Range("A1:AB6").Select() | Out-Null

# --- 2. On the TOR sheet, clear out the two replaced players (rows 5-6),
#        leaving the rows in place but empty. ---

$tor.Rows("5:6").ClearContents() | Out-Null
$tor.Range("A5:XFD6").Select() | Out-Null

# --- 3. On the Non-TOR sheet, delete the two players that were swapped
#        out (rows 2-3) entirely, shifting the remaining players up. ---

$nontor.Rows("2:3").Delete() | Out-Null
$nontor.Activate()
$nontor.Range("A2:XFD3").Select() | Out-Null
